$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 786
$ws1.Range("F6").Value = 17
$ws1.Range("F7").Value = 164
$ws1.Range("F8").Value = 349
$ws1.Range("F10").Value = 518
$ws1.Range("F12").Value = 11875
$ws1.Range("F13").Value = 5426

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 107

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 786
$ws4.Range("F4").Value = 107
$ws4.Range("F8").Value = 17
$ws4.Range("F9").Value = 164
$ws4.Range("F10").Value = 349
$ws4.Range("F12").Value = 518
$ws4.Range("F14").Value = 11875
$ws4.Range("F16").Value = 5426
